$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1337.25
$ws.Range("I40").Value = 1337.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1337.25
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = -1162.25
$ws.Range("N40").ClearContents()
$ws.Range("H51").Value = 5497.0557
$ws.Range("I51").Value = 4979.146
$ws.Range("J51").Value = 6532.875
$ws.Range("K51").Value = 4979.146
$ws.Range("L51").Value = 6532.875
$ws.Range("M51").Value = -4495.146
$ws.Range("N51").Value = -7500.875
$ws.Range("H97").Value = 7900
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H138").Value = 2481.2903
$ws.Range("I138").Value = 2185.5
$ws.Range("J138").Value = 2622.1428
$ws.Range("K138").Value = 6556.5
$ws.Range("L138").Value = 7866.428400000001
$ws.Range("M138").Value = -1416.5
$ws.Range("N138").Value = -18146.4284

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 141626.86
$ws.Range("J32").Value = 21500
$ws.Range("L32").Value = 21500
$ws.Range("N32").Value = -22074
$ws.Range("H34").Value = 184750
$ws.Range("I34").Value = 184750
$ws.Range("K34").Value = 184750
$ws.Range("M34").Value = -184479
$ws.Range("H45").Value = 1546.8948
$ws.Range("I45").Value = 1337.1666
$ws.Range("K45").Value = 1337.1666
$ws.Range("M45").Value = -960.1666
$ws.Range("H61").Value = 4172282
$ws.Range("I61").Value = 6110.684
$ws.Range("K61").Value = 6110.684
$ws.Range("M61").Value = -5898.684
$ws.Range("H69").Value = 200000
$ws.Range("J69").Value = 200000
$ws.Range("L69").Value = 200000
$ws.Range("N69").Value = -201498
$ws.Range("H72").Value = 200000
$ws.Range("J72").Value = 200000
$ws.Range("L72").Value = 600000
$ws.Range("N72").Value = -607488
$ws.Range("H136").Value = 4172282
$ws.Range("I136").Value = 6110.684
$ws.Range("K136").Value = 18332.052
$ws.Range("M136").Value = -15782.052

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 4999.125
$ws.Range("J46").Value = 4999.125
$ws.Range("L46").Value = 4999.125
$ws.Range("N46").Value = -5595.125

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 8469.429
$ws.Range("I105").Value = 9830.637000000001
$ws.Range("K105").Value = 9830.637000000001
$ws.Range("M105").Value = -8083.637000000001
$ws.Range("H107").Value = 880.3333
$ws.Range("I107").Value = 880.3333
$ws.Range("K107").Value = 880.3333
$ws.Range("M107").Value = 1039.6667
$ws.Range("H140").Value = 88298
$ws.Range("J140").Value = 88298
$ws.Range("L140").Value = 88298
$ws.Range("N140").Value = -98658
$ws.Range("H141").Value = 313479.6
$ws.Range("J141").Value = 366849.62
$ws.Range("L141").Value = 366849.62
$ws.Range("N141").Value = -377209.62

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13440549
$ws.Range("I4").Value = 23957802
$ws.Range("K4").Value = 71873406
$ws.Range("M4").Value = -71873294
$ws.Range("H5").Value = 3470410
$ws.Range("I5").Value = 3247317.8
$ws.Range("K5").Value = 9741953.399999999
$ws.Range("M5").Value = -9741841.399999999
$ws.Range("H132").Value = 14487.5
$ws.Range("I132").Value = 1300
$ws.Range("J132").Value = 22400
$ws.Range("K132").Value = 11700
$ws.Range("L132").Value = 201600
$ws.Range("M132").Value = -9170
$ws.Range("N132").Value = -206660
$ws.Range("H135").Value = 3470410
$ws.Range("I135").Value = 3247317.8
$ws.Range("K135").Value = 29225860.2
$ws.Range("M135").Value = -29223325.2

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 649.4
$ws.Range("I97").Value = 718.125
$ws.Range("J97").Value = 374.5
$ws.Range("K97").Value = 718.125
$ws.Range("L97").Value = 374.5
$ws.Range("M97").Value = -222.125
$ws.Range("N97").Value = -1366.5
$ws.Range("H102").Value = 1896.5
$ws.Range("I102").Value = 1852.36
$ws.Range("K102").Value = 1852.36
$ws.Range("M102").Value = -230.3599999999999
$ws.Range("H122").Value = 11538
$ws.Range("I122").Value = 12946.444
$ws.Range("J122").Value = 5200
$ws.Range("K122").Value = 38839.33199999999
$ws.Range("L122").Value = 15600
$ws.Range("M122").Value = -36389.33199999999
$ws.Range("N122").Value = -20500
$ws.Range("H126").Value = 7953.769
$ws.Range("I126").Value = 8809
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 26427
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -23957
$ws.Range("N126").Value = -14690
$ws.Range("H132").Value = 8964.473
$ws.Range("I132").Value = 6500.6665
$ws.Range("K132").Value = 19501.9995
$ws.Range("M132").Value = -16971.9995

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1471.6666
$ws.Range("I30").Value = 1471.6666
$ws.Range("K30").Value = 1471.6666
$ws.Range("M30").Value = -1363.6666
$ws.Range("H35").Value = 18708.5
$ws.Range("I35").Value = 8266.333000000001
$ws.Range("K35").Value = 8266.333000000001
$ws.Range("M35").Value = -7930.333000000001
$ws.Range("H46").Value = 2652.1667
$ws.Range("I46").Value = 561.6667
$ws.Range("J46").Value = 6833.1665
$ws.Range("K46").Value = 561.6667
$ws.Range("L46").Value = 6833.1665
$ws.Range("M46").Value = -373.6667
$ws.Range("N46").Value = -7209.1665
$ws.Range("H59").Value = 73650
$ws.Range("J59").Value = 73650
$ws.Range("L59").Value = 73650
$ws.Range("N59").Value = -74958
$ws.Range("H61").Value = 22843.857
$ws.Range("I61").Value = 22843.857
$ws.Range("K61").Value = 22843.857
$ws.Range("M61").Value = -22641.857
$ws.Range("H68").Value = 2064.4
$ws.Range("I68").Value = 1749.8
$ws.Range("J68").Value = 2379
$ws.Range("K68").Value = 1749.8
$ws.Range("L68").Value = 2379
$ws.Range("M68").Value = -1000.8
$ws.Range("N68").Value = -3877
$ws.Range("H71").Value = 2064.4
$ws.Range("I71").Value = 1749.8
$ws.Range("J71").Value = 2379
$ws.Range("K71").Value = 8749
$ws.Range("L71").Value = 11895
$ws.Range("M71").Value = -5005
$ws.Range("N71").Value = -19383
$ws.Range("H93").Value = 1811.5161
$ws.Range("I93").Value = 1484.9231
$ws.Range("J93").Value = 3509.8
$ws.Range("K93").Value = 1484.9231
$ws.Range("L93").Value = 3509.8
$ws.Range("M93").Value = -236.9231
$ws.Range("N93").Value = -6005.8
$ws.Range("H100").Value = 1560.4
$ws.Range("I100").Value = 1551.5834
$ws.Range("J100").Value = 1595.6666
$ws.Range("K100").Value = 1551.5834
$ws.Range("L100").Value = 1595.6666
$ws.Range("M100").Value = -1010.5834
$ws.Range("N100").Value = -2677.6666
$ws.Range("H113").Value = 22843.857
$ws.Range("I113").Value = 22843.857
$ws.Range("K113").Value = 22843.857
$ws.Range("M113").Value = -20673.857
$ws.Range("H122").Value = 3500.6924
$ws.Range("I122").Value = 3152
$ws.Range("J122").Value = 4663
$ws.Range("K122").Value = 9456
$ws.Range("L122").Value = 13989
$ws.Range("M122").Value = -7006
$ws.Range("N122").Value = -18889

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 54960.332
$ws.Range("J125").Value = 54999.145
$ws.Range("L125").Value = 54999.145
$ws.Range("N125").Value = -64839.145
$ws.Range("H126").Value = 3235.7856
$ws.Range("I126").Value = 3681.889
$ws.Range("J126").Value = 2432.8
$ws.Range("K126").Value = 11045.667
$ws.Range("L126").Value = 7298.400000000001
$ws.Range("M126").Value = -8575.667000000001
$ws.Range("N126").Value = -12238.4
